# This script rotates the list of account numbers in column A (rows 2-141)
# of the active worksheet. The first 28 values (rows 2-29) are moved to the
# bottom of the list (becoming the last 28 rows), while the remaining values
# (rows 30-141) shift up to the top. Net effect: the values that formerly
# started at row 30 now start at row 2, and the old top-28 block is appended
# at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 141
$rotateCount = 28   # number of rows to move from the top to the bottom
$totalRows = $lastDataRow - $firstDataRow + 1

# Read all current values from A2:A141 into an array (in row order).
$values = @()
for ($i = 0; $i -lt $totalRows; $i++) {
    $values += $ws.Cells.Item($firstDataRow + $i, 1).Value2
}

# Build the rotated array: drop the first $rotateCount items from the front
# and append them to the end.
$rotated = @()
for ($i = 0; $i -lt $totalRows; $i++) {
    $srcIndex = ($i + $rotateCount) % $totalRows
    $rotated += $values[$srcIndex]
}

# Write the rotated values back into the same range, one cell at a time.
for ($i = 0; $i -lt $totalRows; $i++) {
    $ws.Cells.Item($firstDataRow + $i, 1).Value = $rotated[$i]
}
